{"js": "// Apply the LOM3230.docx edits described by the commit diff.\n\nconst body = context.document.body;\n\n// Helper: find the index (within body.paragraphs) of the paragraph whose\n// text contains the given marker substring.\nasync function findParagraphIndex(marker) {\n  const paras = body.paragraphs;\n  paras.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < paras.items.length; i++) {\n    if (paras.items[i].text.indexOf(marker) !== -1) {\n      return { paras, index: i };\n    }\n  }\n  return { paras, index: -1 };\n}\n\n// --- 1. \"Ativa\u00e7\u00e3o: 01/01/2012\" -> \"Ativa\u00e7\u00e3o: 01/01/2023\" ---\nconst ativacao = body.search(\"Ativa\u00e7\u00e3o: 01/01/2012\", { matchCase: true });\nativacao.load(\"text\");\nawait context.sync();\nif (ativacao.items.length > 0) {\n  ativacao.items[0].insertText(\"Ativa\u00e7\u00e3o: 01/01/2023\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 2. Append an italic EN objectives paragraph after \"Objetivos\" text paragraph ---\nconst objEnText =\n  \"To present experimental techniques for the characterization of electrical, magnetic and thermal  properties of materials.\";\n\n{\n  const { paras, index } = await findParagraphIndex(\n    \"Apresentar as t\u00e9cnicas experimentais de caracteriza\u00e7\u00e3o de propriedades el\u00e9tricas, magn\u00e9ticas, t\u00e9rmicas e \u00f3pticas de materiais.\"\n  );\n  if (index !== -1) {\n    const p = paras.items[index].insertParagraph(objEnText, Word.InsertLocation.after);\n    p.font.italic = true;\n    await context.sync();\n  }\n}\n\n// --- 3. Append the same italic EN paragraph after \"Programa resumido\" text paragraph ---\n{\n  const { paras, index } = await findParagraphIndex(\n    \"Estudo das t\u00e9cnicas de caracteriza\u00e7\u00e3o de propriedades el\u00e9tricas, magn\u00e9ticas, t\u00e9rmicas e \u00f3pticas de materiais.\"\n  );\n  if (index !== -1) {\n    const p = paras.items[index].insertParagraph(objEnText, Word.InsertLocation.after);\n    p.font.italic = true;\n    await context.sync();\n  }\n}\n\n// --- 4. Replace the \"Programa\" paragraph body and append an italic EN paragraph ---\n{\n  const programaNewText =\n    \"Propriedades el\u00e9tricas: condutividade el\u00e9trica em metais puros, ligas met\u00e1licas e semicondutores,  e supercondutores; Efeito Hall; Lei de Ohm e depend\u00eancia com a temperatura.Propriedades magn\u00e9ticas: susceptibilidade magn\u00e9tica e magnetiza\u00e7\u00e3o c.c. Curvas de histerese de materiais magn\u00e9ticos macios. Medidas de magnetostric\u00e7\u00e3o.Propriedades t\u00e9rmicas dos materiais:  expans\u00e3o t\u00e9rmica.\";\n  const programaEnText =\n    \"histerese de materiais magn\u00e9ticos macios. Medidas de magnetostric\u00e7\u00e3o.Propriedades t\u00e9rmicas dos materiais:  expans\u00e3o t\u00e9rmica.Electrical properties: electrical conductivity in pure metals, metallic alloys and semiconductors, and superconductors; Hall Effect; Ohm's Law and dependence on temperature.Magnetic properties: magnetic susceptibility and c.c. magnetization. Hysteresis curves of soft magnetic materials. Magnetostriction measurements.Thermal properties of materials: thermal expansion.\";\n\n  const { paras, index } = await findParagraphIndex(\n    \"Propriedades el\u00e9tricas: condutividade el\u00e9trica em metais puros, ligas met\u00e1licas, semicondutores, isolantes e supercondutores\"\n  );\n  if (index !== -1) {\n    const target = paras.items[index];\n    target.insertText(programaNewText, Word.InsertLocation.replace);\n    const p = target.insertParagraph(programaEnText, Word.InsertLocation.after);\n    p.font.italic = true;\n    await context.sync();\n  }\n}\n\n// --- 5. \"Crit\u00e9rio\" text ---\nconst criterioOld = body.search(\n  \"M\u00e9dia aritm\u00e9tica de duas provas escritas, testes, trabalhos e relat\u00f3rios: P1, P2 e TR. Conceito Final = (P1 + P2 + TR)/3\",\n  { matchCase: true }\n);\ncriterioOld.load(\"text\");\nawait context.sync();\nif (criterioOld.items.length > 0) {\n  criterioOld.items[0].insertText(\n    \"M\u00e9dia aritm\u00e9tica das notas dos relat\u00f3rios de cada experimento\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// --- 6. Bibliografia paragraph: collapse the 8 bullet-like lines into one run ---\n{\n  const bibliografiaNewText =\n    \"HUMMEL, R. E. Electronic Properties of Materials, Springer, 2000.RAYMOND A. SERWAY, CLEMENT J. MOSES, CURT A. MOYER. Modern Physics 3rd Edition,  Cengage Learning, Inc., 2005.SOLYMAR, L.; WALSH, D. Electrical Properties of Materials, Oxford University Press, 2009.NICOLA A. SPALDIN, Magnetic Materials, Fundamentals and Applications, SECOND EDITION, Cambridge University Press, 2011ROBERT, P. Electrical and Magnetic Properties of Materials, Artech House, 1998.SPEYER, R. Thermal Analysis of Materials, CRC Press, 1993.\";\n\n  const { paras, index } = await findParagraphIndex(\n    \"HUMMEL, R. E. Electronic Properties of Materials, Springer, 2000.\"\n  );\n  if (index !== -1) {\n    paras.items[index].insertText(bibliografiaNewText, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# Apply the LOM3230.docx edits described by the commit diff.\n\n$d = $word.ActiveDocument\n\nfunction FindParaIndex($marker) {\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($p.Range.Text.Contains($marker)) {\n            return $i\n        }\n    }\n    return -1\n}\n\n# --- 1. \"Ativa\u00e7\u00e3o: 01/01/2012\" -> \"Ativa\u00e7\u00e3o: 01/01/2023\" ---\n$range = $d.Content\n$range.Find.Execute(\"Ativa\u00e7\u00e3o: 01/01/2012\", $false, $false, $false, $false, $false, $true, 1, $false, \"Ativa\u00e7\u00e3o: 01/01/2023\", 2) | Out-Null\n\n# --- 2. Append an italic EN objectives paragraph after the \"Objetivos\" text paragraph ---\n$objEnText = \"To present experimental techniques for the characterization of electrical, magnetic and thermal  properties of materials.\"\n\n$idx = FindParaIndex(\"Apresentar as t\u00e9cnicas experimentais de caracteriza\u00e7\u00e3o de propriedades el\u00e9tricas, magn\u00e9ticas, t\u00e9rmicas e \u00f3pticas de materiais.\")\nif ($idx -ne -1) {\n    $p = $d.Paragraphs.Item($idx)\n    $p.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs.Item($idx + 1)\n    $newPara.Range.Text = $objEnText\n    $textOnly = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)\n    $textOnly.Font.Italic = $true\n}\n\n# --- 3. Append the same italic EN paragraph after the \"Programa resumido\" text paragraph ---\n$idx = FindParaIndex(\"Estudo das t\u00e9cnicas de caracteriza\u00e7\u00e3o de propriedades el\u00e9tricas, magn\u00e9ticas, t\u00e9rmicas e \u00f3pticas de materiais.\")\nif ($idx -ne -1) {\n    $p = $d.Paragraphs.Item($idx)\n    $p.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs.Item($idx + 1)\n    $newPara.Range.Text = $objEnText\n    $textOnly = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)\n    $textOnly.Font.Italic = $true\n}\n\n# --- 4. Replace the \"Programa\" paragraph body and append an italic EN paragraph ---\n$programaNewText = \"Propriedades el\u00e9tricas: condutividade el\u00e9trica em metais puros, ligas met\u00e1licas e semicondutores,  e supercondutores; Efeito Hall; Lei de Ohm e depend\u00eancia com a temperatura.Propriedades magn\u00e9ticas: susceptibilidade magn\u00e9tica e magnetiza\u00e7\u00e3o c.c. Curvas de histerese de materiais magn\u00e9ticos macios. Medidas de magnetostric\u00e7\u00e3o.Propriedades t\u00e9rmicas dos materiais:  expans\u00e3o t\u00e9rmica.\"\n$programaEnText = \"histerese de materiais magn\u00e9ticos macios. Medidas de magnetostric\u00e7\u00e3o.Propriedades t\u00e9rmicas dos materiais:  expans\u00e3o t\u00e9rmica.Electrical properties: electrical conductivity in pure metals, metallic alloys and semiconductors, and superconductors; Hall Effect; Ohm's Law and dependence on temperature.Magnetic properties: magnetic susceptibility and c.c. magnetization. Hysteresis curves of soft magnetic materials. Magnetostriction measurements.Thermal properties of materials: thermal expansion.\"\n\n$idx = FindParaIndex(\"Propriedades el\u00e9tricas: condutividade el\u00e9trica em metais puros, ligas met\u00e1licas, semicondutores, isolantes e supercondutores\")\nif ($idx -ne -1) {\n    $p = $d.Paragraphs.Item($idx)\n    $p.Range.Text = $programaNewText\n    $p.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs.Item($idx + 1)\n    $newPara.Range.Text = $programaEnText\n    $textOnly = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)\n    $textOnly.Font.Italic = $true\n}\n\n# --- 5. \"Crit\u00e9rio\" text ---\n$range = $d.Content\n$range.Find.Execute(\"M\u00e9dia aritm\u00e9tica de duas provas escritas, testes, trabalhos e relat\u00f3rios: P1, P2 e TR. Conceito Final = (P1 + P2 + TR)/3\", $false, $false, $false, $false, $false, $true, 1, $false, \"M\u00e9dia aritm\u00e9tica das notas dos relat\u00f3rios de cada experimento\", 2) | Out-Null\n\n# --- 6. Bibliografia paragraph: collapse the 8 bullet-like lines into one run ---\n$bibliografiaNewText = \"HUMMEL, R. E. Electronic Properties of Materials, Springer, 2000.RAYMOND A. SERWAY, CLEMENT J. MOSES, CURT A. MOYER. Modern Physics 3rd Edition,  Cengage Learning, Inc., 2005.SOLYMAR, L.; WALSH, D. Electrical Properties of Materials, Oxford University Press, 2009.NICOLA A. SPALDIN, Magnetic Materials, Fundamentals and Applications, SECOND EDITION, Cambridge University Press, 2011ROBERT, P. Electrical and Magnetic Properties of Materials, Artech House, 1998.SPEYER, R. Thermal Analysis of Materials, CRC Press, 1993.\"\n\n$idx = FindParaIndex(\"HUMMEL, R. E. Electronic Properties of Materials, Springer, 2000.\")\nif ($idx -ne -1) {\n    $p = $d.Paragraphs.Item($idx)\n    $p.Range.Text = $bibliografiaNewText\n}\n"}
